# edit.ps1 - Append 45 new data rows (102-146) to the
# "master-reg_center_machine_devic" sheet, matching the target diff:
#   - New rows 102-146, columns A-H:
#       A: cycles 10002..10010
#       B: cycles 10021..10029
#       C: sequential 3000121..3000165
#       D: "eng" (shared string, same as existing rows)
#       E: TRUE (boolean)
#       F: "superadmin()" (shared string)
#       G/H: "now()" (shared string)
#   - dimension grows to A1:H146 (automatic from writing the cells)
#   - sheetView selection becomes A147:XFD1048576 (as if the user
#     clicked/selected the row right after the last data row)
#   - pageSetup orientation explicitly set to portrait

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRows = @(
    @(10002, 10021, 3000121),
    @(10003, 10022, 3000122),
    @(10004, 10023, 3000123),
    @(10005, 10024, 3000124),
    @(10006, 10025, 3000125),
    @(10007, 10026, 3000126),
    @(10008, 10027, 3000127),
    @(10009, 10028, 3000128),
    @(10010, 10029, 3000129),
    @(10002, 10021, 3000130),
    @(10003, 10022, 3000131),
    @(10004, 10023, 3000132),
    @(10005, 10024, 3000133),
    @(10006, 10025, 3000134),
    @(10007, 10026, 3000135),
    @(10008, 10027, 3000136),
    @(10009, 10028, 3000137),
    @(10010, 10029, 3000138),
    @(10002, 10021, 3000139),
    @(10003, 10022, 3000140),
    @(10004, 10023, 3000141),
    @(10005, 10024, 3000142),
    @(10006, 10025, 3000143),
    @(10007, 10026, 3000144),
    @(10008, 10027, 3000145),
    @(10009, 10028, 3000146),
    @(10010, 10029, 3000147),
    @(10002, 10021, 3000148),
    @(10003, 10022, 3000149),
    @(10004, 10023, 3000150),
    @(10005, 10024, 3000151),
    @(10006, 10025, 3000152),
    @(10007, 10026, 3000153),
    @(10008, 10027, 3000154),
    @(10009, 10028, 3000155),
    @(10010, 10029, 3000156),
    @(10002, 10021, 3000157),
    @(10003, 10022, 3000158),
    @(10004, 10023, 3000159),
    @(10005, 10024, 3000160),
    @(10006, 10025, 3000161),
    @(10007, 10026, 3000162),
    @(10008, 10027, 3000163),
    @(10009, 10028, 3000164),
    @(10010, 10029, 3000165)
)

$startRow = 102
for ($i = 0; $i -lt $dataRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $dataRows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]          # A: regcntr_id
    $ws.Cells.Item($r, 2).Value = $vals[1]          # B: machine_id
    $ws.Cells.Item($r, 3).Value = $vals[2]          # C: device_id
    $ws.Cells.Item($r, 4).Value = "eng"             # D: lang_code
    $ws.Cells.Item($r, 5).Value = $true              # E: is_active
    $ws.Cells.Item($r, 6).Value = "superadmin()"     # F: cr_by
    $ws.Cells.Item($r, 7).Value = "now()"            # G: cr_dtimes
    $ws.Cells.Item($r, 8).Value = "now()"            # H: eff_dtimes
}

# Mirror the post-entry selection state captured in the diff: the cursor
# lands on the row after the last data row, selected down to the bottom
# of the sheet.
[void]$ws.Range("A147:XFD1048576").Select()

# The diff also flips the sheet's page orientation to portrait.
$ws.PageSetup.Orientation = 1

